$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.370.70"
$ws.Range("E2").Value = "'  -2.66%  "
$ws.Range("D3").Value = "'3.546.02"
$ws.Range("E3").Value = "'  -3.64%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'582.22"
$ws.Range("E5").Value = "'  -0.26%  "
$ws.Range("D6").Value = "'173.15"
$ws.Range("E6").Value = "'  -2.31%  "
$ws.Range("D7").Value = "'3.538.50"
$ws.Range("E7").Value = "'  -3.60%  "
$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "'  -1.04%  "
$ws.Range("E9").Value = "'  -0.01%  "
$ws.Range("E10").Value = "'  -4.87%  "
$ws.Range("D11").Value = "'6.76"
$ws.Range("E11").Value = "'  -2.10%  "
$ws.Range("D12").Value = "'0.589"
$ws.Range("E12").Value = "'  -3.51%  "
$ws.Range("D13").Value = "'47.84"
$ws.Range("E13").Value = "'  -2.30%  "
$ws.Range("E14").Value = "'  -5.65%  "
$ws.Range("D15").Value = "'4.115.06"
$ws.Range("E15").Value = "'  -3.69%  "
$ws.Range("D16").Value = "'8.59"
$ws.Range("E16").Value = "'  -4.84%  "
$ws.Range("D17").Value = "'631.01"
$ws.Range("E17").Value = "'  -6.60%  "
$ws.Range("B18").Value = "'WrappedBTC"
$ws.Range("C18").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'69.462.86"
$ws.Range("E18").Value = "'  -2.55%  "
$ws.Range("B19").Value = "'WrappedEther"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.536.13"
$ws.Range("E19").Value = "'  -4.07%  "
$ws.Range("E20").Value = "'  +0.22%  "
$ws.Range("D21").Value = "'17.52"
$ws.Range("E21").Value = "'  -2.26%  "
$ws.Range("D22").Value = "'11.30"
$ws.Range("E22").Value = "'  -2.25%  "
$ws.Range("D23").Value = "'0.894"
$ws.Range("E23").Value = "'  -4.86%  "
$ws.Range("D24").Value = "'16.08"
$ws.Range("E24").Value = "'  -6.97%  "
$ws.Range("D25").Value = "'98.01"
$ws.Range("E25").Value = "'  -3.65%  "
$ws.Range("E26").Value = "'  -3.53%  "
$ws.Range("E27").Value = "'  +0.11%  "
$ws.Range("E28").Value = "'  -6.11%  "
$ws.Range("E29").Value = "'  -7.44%  "
$ws.Range("D30").Value = "'33.02"
$ws.Range("E30").Value = "'  -5.87%  "
$ws.Range("D31").Value = "'3.18"
$ws.Range("E31").Value = "'  -6.74%  "
$ws.Range("D32").Value = "'8.61"
$ws.Range("E32").Value = "'  -5.53%  "
$ws.Range("E33").Value = "'  -6.30%  "
$ws.Range("E34").Value = "'  -6.68%  "
$ws.Range("D35").Value = "'640.67"
$ws.Range("E35").Value = "'  +9.63%  "
$ws.Range("D36").Value = "'10.84"
$ws.Range("E36").Value = "'  -3.16%  "
$ws.Range("D37").Value = "'3.53"
$ws.Range("E37").Value = "'  -11.55%  "
$ws.Range("E38").Value = "'  -4.25%  "
$ws.Range("D39").Value = "'57.47"
$ws.Range("E39").Value = "'  -2.10%  "
$ws.Range("E40").Value = "'  +0.10%  "
$ws.Range("E41").Value = "'  -1.39%  "
$ws.Range("D42").Value = "'0.136"
$ws.Range("E42").Value = "'  -4.49%  "
$ws.Range("D43").Value = "'3.410.40"
$ws.Range("E43").Value = "'  -5.80%  "
$ws.Range("D45").Value = "'32.98"
$ws.Range("E45").Value = "'  -6.20%  "
$ws.Range("E46").Value = "'  -7.28%  "
$ws.Range("E47").Value = "'  -6.27%  "
$ws.Range("E48").Value = "'  -4.71%  "
$ws.Range("D49").Value = "'0.131"
$ws.Range("E49").Value = "'  -2.64%  "
$ws.Range("D50").Value = "'133.21"
$ws.Range("E50").Value = "'  -1.07%  "
$ws.Range("D51").Value = "'5.65"
$ws.Range("E51").Value = "'  +13.13%  "
